# Fruta / hortaliza, semanal
# Insert one new weekly record at row 26 (pushing the existing rows 26..103
# down to 27..104) on the single data sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26; Excel shifts rows 26..103 down to
# 27..104 and grows the used range to A1:R104 automatically.
$ws.Rows("26:26").Insert()

# Reuse the values already sitting in the row that got pushed down to 27
# (same market/region/category/variety/quality/unit-origin/classification)
# so we don't have to retype the accented Spanish text, then overwrite the
# few cells that actually hold new data for this new record.
$ws.Range("A26").Value = $ws.Range("A27").Value2
$ws.Range("B26").Value = $ws.Range("B27").Value2
$ws.Range("C26").Value = $ws.Range("C27").Value2
$ws.Range("E26").Value = $ws.Range("E27").Value2
$ws.Range("F26").Value = $ws.Range("F27").Value2
$ws.Range("G26").Value = $ws.Range("G27").Value2
$ws.Range("H26").Value = $ws.Range("H27").Value2
$ws.Range("I26").Value = $ws.Range("I27").Value2
$ws.Range("O26").Value = $ws.Range("O27").Value2
$ws.Range("R26").Value = $ws.Range("R27").Value2

# New record's own data.
$ws.Range("D26").Value = 44560
$ws.Range("J26").Value = 720
$ws.Range("K26").Value = 11000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 11500
$ws.Range("N26").Value = "$/malla 70 unidades"
$ws.Range("P26").Value = 164
$ws.Range("Q26").Value = 70
